$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; remove protection temporarily so the values can be
# updated, then restore it.
$ws.Unprotect()

# Update the confidential disclosure date text (cell A9)
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-06 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-6
$ws.Range("D2").Value = 0.2552645641516195
$ws.Range("E2").Value = 0.008460825220213186

$ws.Range("D3").Value = 0.2528416936818407
$ws.Range("E3").Value = 0.01485282203618676

$ws.Range("D4").Value = 0.2455637014722346
$ws.Range("E4").Value = 0.00695450594030711

$ws.Range("D5").Value = 0.2463300406943051
$ws.Range("E5").Value = 0.01172027607761428

$ws.Range("E6").Value = 0.01050999184554491

$ws.Protect()
